$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 15.0597367646986
$ws.Cells.Item(2, 3).Value = 12.03898504284022
$ws.Cells.Item(2, 4).Value = 4.864864825314219
$ws.Cells.Item(2, 6).Value = 25.21674136324056
$ws.Cells.Item(2, 7).Value = 3.631373362726689
$ws.Cells.Item(2, 12).Value = 10.8403660579065
$ws.Cells.Item(2, 13).Value = 15.26941240228047
$ws.Cells.Item(2, 14).Value = 18.17974642313193
$ws.Cells.Item(2, 15).Value = 22.4545918845222
$ws.Cells.Item(3, 2).Value = 14.61073592849474
$ws.Cells.Item(3, 3).Value = 11.90198604048744
$ws.Cells.Item(3, 4).Value = 4.835071270709284
$ws.Cells.Item(3, 6).Value = 25.17937328717281
$ws.Cells.Item(3, 7).Value = 3.633526337783374
$ws.Cells.Item(3, 12).Value = 10.85543088089858
$ws.Cells.Item(3, 13).Value = 15.18983167809014
$ws.Cells.Item(3, 14).Value = 18.23925466731086
$ws.Cells.Item(3, 15).Value = 22.48466969989769
$ws.Cells.Item(4, 2).Value = 14.33040571121336
$ws.Cells.Item(4, 3).Value = 11.81639834605854
$ws.Cells.Item(4, 4).Value = 4.816442859591607
$ws.Cells.Item(4, 6).Value = 25.16403354563555
$ws.Cells.Item(4, 7).Value = 3.634919241950064
$ws.Cells.Item(4, 12).Value = 10.86629914729202
$ws.Cells.Item(4, 13).Value = 15.14329069667549
$ws.Cells.Item(4, 14).Value = 18.27765409218454
$ws.Cells.Item(4, 15).Value = 22.50901248577385
$ws.Cells.Item(5, 2).Value = 14.21518125954959
$ws.Cells.Item(5, 3).Value = 11.78117228687878
$ws.Cells.Item(5, 4).Value = 4.80876989485342
$ws.Cells.Item(5, 6).Value = 25.15969871212946
$ws.Cells.Item(5, 7).Value = 3.635504764226553
$ws.Cells.Item(5, 12).Value = 10.87113515742496
$ws.Cells.Item(5, 13).Value = 15.12492342910498
$ws.Cells.Item(5, 14).Value = 18.29377138337781
$ws.Cells.Item(5, 15).Value = 22.52040627899009
$ws.Cells.Item(6, 2).Value = 14.19599416313503
$ws.Cells.Item(6, 3).Value = 11.77530261117033
$ws.Cells.Item(6, 4).Value = 4.807490941824442
$ws.Cells.Item(6, 6).Value = 25.15909473207471
$ws.Cells.Item(6, 7).Value = 3.635603072694467
$ws.Cells.Item(6, 12).Value = 10.87196276426623
$ws.Cells.Item(6, 13).Value = 15.12191012386907
$ws.Cells.Item(6, 14).Value = 18.29647602270961
$ws.Cells.Item(6, 15).Value = 22.5223871192584
$ws.Cells.Item(7, 2).Value = 14.32885550655602
$ws.Cells.Item(7, 3).Value = 11.81592465743346
$ws.Cells.Item(7, 4).Value = 4.816339706935453
$ws.Cells.Item(7, 6).Value = 25.16396732192849
$ws.Cells.Item(7, 7).Value = 3.634927065946082
$ws.Cells.Item(7, 12).Value = 10.8663627189239
$ws.Cells.Item(7, 13).Value = 15.14304054691245
$ws.Cells.Item(7, 14).Value = 18.27786955416461
$ws.Cells.Item(7, 15).Value = 22.50916018360189
$ws.Cells.Item(8, 2).Value = 14.90598868316775
$ws.Cells.Item(8, 3).Value = 11.99206544227083
$ws.Cells.Item(8, 4).Value = 4.854663050205708
$ws.Cells.Item(8, 6).Value = 25.20228156800196
$ws.Cells.Item(8, 7).Value = 3.632101013222338
$ws.Cells.Item(8, 12).Value = 10.84522461060928
$ws.Cells.Item(8, 13).Value = 15.24150010772495
$ws.Cells.Item(8, 14).Value = 18.19987932407169
$ws.Cells.Item(8, 15).Value = 22.46374153535412
$ws.Cells.Item(9, 2).Value = 15.99365395850863
$ws.Cells.Item(9, 3).Value = 12.32490546231159
$ws.Cells.Item(9, 4).Value = 4.927047916343796
$ws.Cells.Item(9, 6).Value = 25.33748736980197
$ws.Cells.Item(9, 7).Value = 3.627119646369504
$ws.Cells.Item(9, 12).Value = 10.81660645488696
$ws.Cells.Item(9, 13).Value = 15.45232252170519
$ws.Cells.Item(9, 14).Value = 18.06165200393496
$ws.Cells.Item(9, 15).Value = 22.42141522945023
$ws.Cells.Item(10, 2).Value = 16.75707238076473
$ws.Cells.Item(10, 3).Value = 12.56059548100226
$ws.Cells.Item(10, 4).Value = 4.97839943285877
$ws.Cells.Item(10, 6).Value = 25.47296018209811
$ws.Cells.Item(10, 7).Value = 3.623797916356144
$ws.Cells.Item(10, 12).Value = 10.80339183400204
$ws.Cells.Item(10, 13).Value = 15.61707495501342
$ws.Cells.Item(10, 14).Value = 17.96898622477031
$ws.Cells.Item(10, 15).Value = 22.4189436407998
$ws.Cells.Item(11, 2).Value = 17.09501943952048
$ws.Cells.Item(11, 3).Value = 12.66565496659931
$ws.Cells.Item(11, 4).Value = 5.001333340544482
$ws.Cells.Item(11, 6).Value = 25.54229267994371
$ws.Cells.Item(11, 7).Value = 3.622359415628085
$ws.Cells.Item(11, 12).Value = 10.79907232159986
$ws.Cells.Item(11, 13).Value = 15.69395029336248
$ws.Cells.Item(11, 14).Value = 17.92874396925298
$ws.Cells.Item(11, 15).Value = 22.42404877747623
$ws.Cells.Item(12, 2).Value = 17.22153230248352
$ws.Cells.Item(12, 3).Value = 12.70510855640318
$ws.Cells.Item(12, 4).Value = 5.009953992999166
$ws.Cells.Item(12, 6).Value = 25.56963982996009
$ws.Cells.Item(12, 7).Value = 3.62182506995144
$ws.Cells.Item(12, 12).Value = 10.79767945477678
$ws.Cells.Item(12, 13).Value = 15.72331932417122
$ws.Cells.Item(12, 14).Value = 17.9137789937466
$ws.Cells.Item(12, 15).Value = 22.42687758132238
$ws.Cells.Item(13, 2).Value = 17.1943521921644
$ws.Cells.Item(13, 3).Value = 12.69662651466892
$ws.Cells.Item(13, 4).Value = 5.008100265851166
$ws.Cells.Item(13, 6).Value = 25.56370181317915
$ws.Cells.Item(13, 7).Value = 3.621939689879879
$ws.Cells.Item(13, 12).Value = 10.79796864124729
$ws.Cells.Item(13, 13).Value = 15.71698301476354
$ws.Cells.Item(13, 14).Value = 17.91698980371214
$ws.Cells.Item(13, 15).Value = 22.42622852700501
$ws.Cells.Item(14, 2).Value = 17.10545762881541
$ws.Cells.Item(14, 3).Value = 12.6689075660517
$ws.Cells.Item(14, 4).Value = 5.002043862890374
$ws.Cells.Item(14, 6).Value = 25.54452074817143
$ws.Cells.Item(14, 7).Value = 3.622315246888367
$ws.Cells.Item(14, 12).Value = 10.79895286561307
$ws.Cells.Item(14, 13).Value = 15.69636144241634
$ws.Cells.Item(14, 14).Value = 17.92750730831965
$ws.Cells.Item(14, 15).Value = 22.42426355955624
$ws.Cells.Item(15, 2).Value = 17.05081366897467
$ws.Cells.Item(15, 3).Value = 12.65188533202696
$ws.Cells.Item(15, 4).Value = 4.998325739001384
$ws.Cells.Item(15, 6).Value = 25.53291356565837
$ws.Cells.Item(15, 7).Value = 3.622546637071248
$ws.Cells.Item(15, 12).Value = 10.79958734153731
$ws.Cells.Item(15, 13).Value = 15.68376315393137
$ws.Cells.Item(15, 14).Value = 17.93398522229603
$ws.Cells.Item(15, 15).Value = 22.4231765749304
$ws.Cells.Item(16, 2).Value = 16.73478853098626
$ws.Cells.Item(16, 3).Value = 12.55368446829343
$ws.Cells.Item(16, 4).Value = 4.976891846379939
$ws.Cells.Item(16, 6).Value = 25.46858290184802
$ws.Cells.Item(16, 7).Value = 3.623893381594883
$ws.Cells.Item(16, 12).Value = 10.80370813239331
$ws.Cells.Item(16, 13).Value = 15.61208813708451
$ws.Cells.Item(16, 14).Value = 17.97165453391948
$ws.Cells.Item(16, 15).Value = 22.41873536221875
$ws.Cells.Item(17, 2).Value = 16.5384384660477
$ws.Cells.Item(17, 3).Value = 12.49287452144023
$ws.Cells.Item(17, 4).Value = 4.963631811492968
$ws.Cells.Item(17, 6).Value = 25.43108114593931
$ws.Cells.Item(17, 7).Value = 3.624738115562117
$ws.Cells.Item(17, 12).Value = 10.80666913831189
$ws.Cells.Item(17, 13).Value = 15.56859818620373
$ws.Cells.Item(17, 14).Value = 17.99525237243283
$ws.Cells.Item(17, 15).Value = 22.4176064470067
$ws.Cells.Item(18, 2).Value = 16.42463062392125
$ws.Cells.Item(18, 3).Value = 12.45769636745046
$ws.Cells.Item(18, 4).Value = 4.955965034617221
$ws.Cells.Item(18, 6).Value = 25.41023743482177
$ws.Cells.Item(18, 7).Value = 3.625230818014741
$ws.Cells.Item(18, 12).Value = 10.80853151056465
$ws.Cells.Item(18, 13).Value = 15.54376663349355
$ws.Cells.Item(18, 14).Value = 18.00900523982374
$ws.Cells.Item(18, 15).Value = 22.41754353225843
$ws.Cells.Item(19, 2).Value = 16.38595118683768
$ws.Cells.Item(19, 3).Value = 12.44575159492637
$ws.Cells.Item(19, 4).Value = 4.953362414225817
$ws.Cells.Item(19, 6).Value = 25.40330529962961
$ws.Cells.Item(19, 7).Value = 3.62539881391408
$ws.Cells.Item(19, 12).Value = 10.8091894468553
$ws.Cells.Item(19, 13).Value = 15.53539106029865
$ws.Cells.Item(19, 14).Value = 18.01369267391299
$ws.Cells.Item(19, 15).Value = 22.41762293975779
$ws.Cells.Item(20, 2).Value = 16.55943145624911
$ws.Cells.Item(20, 3).Value = 12.49936888476619
$ws.Cells.Item(20, 4).Value = 4.965047519300908
$ws.Cells.Item(20, 6).Value = 25.43499820857286
$ws.Cells.Item(20, 7).Value = 3.624647485295074
$ws.Cells.Item(20, 12).Value = 10.80633745233888
$ws.Cells.Item(20, 13).Value = 15.57320899392012
$ws.Cells.Item(20, 14).Value = 17.99272171890482
$ws.Cells.Item(20, 15).Value = 22.41766592991576
$ws.Cells.Item(21, 2).Value = 17.13160863002368
$ws.Cells.Item(21, 3).Value = 12.6770584094163
$ws.Cells.Item(21, 4).Value = 5.003824530928078
$ws.Cells.Item(21, 6).Value = 25.55012517680834
$ws.Cells.Item(21, 7).Value = 3.622204655305663
$ws.Cells.Item(21, 12).Value = 10.79865718845918
$ws.Cells.Item(21, 13).Value = 15.70241165177053
$ws.Cells.Item(21, 14).Value = 17.9244106370041
$ws.Cells.Item(21, 15).Value = 22.42481641715408
$ws.Cells.Item(22, 2).Value = 17.49700114001765
$ws.Cells.Item(22, 3).Value = 12.79125558866853
$ws.Cells.Item(22, 4).Value = 5.028793768603635
$ws.Cells.Item(22, 6).Value = 25.63172634300591
$ws.Cells.Item(22, 7).Value = 3.620668622712638
$ws.Cells.Item(22, 12).Value = 10.79505287944907
$ws.Cells.Item(22, 13).Value = 15.78834770522721
$ws.Cells.Item(22, 14).Value = 17.88136139136106
$ws.Cells.Item(22, 15).Value = 22.4347092791321
$ws.Cells.Item(23, 2).Value = 17.30280322164755
$ws.Cells.Item(23, 3).Value = 12.73048971696543
$ws.Cells.Item(23, 4).Value = 5.015502281760782
$ws.Cells.Item(23, 6).Value = 25.58759803528692
$ws.Cells.Item(23, 7).Value = 3.62148291380998
$ws.Cells.Item(23, 12).Value = 10.79684724603018
$ws.Cells.Item(23, 13).Value = 15.74235173786046
$ws.Cells.Item(23, 14).Value = 17.90419188589667
$ws.Cells.Item(23, 15).Value = 22.4289519406487
$ws.Cells.Item(24, 2).Value = 16.54994339300698
$ws.Cells.Item(24, 3).Value = 12.49643345971665
$ws.Cells.Item(24, 4).Value = 4.964407612793046
$ws.Cells.Item(24, 6).Value = 25.43322507137605
$ws.Cells.Item(24, 7).Value = 3.624688437263416
$ws.Cells.Item(24, 12).Value = 10.80648690899455
$ws.Cells.Item(24, 13).Value = 15.57112391238292
$ws.Cells.Item(24, 14).Value = 17.99386524740355
$ws.Cells.Item(24, 15).Value = 22.41763721195172
$ws.Cells.Item(25, 2).Value = 15.70509241807629
$ws.Cells.Item(25, 3).Value = 12.23633336809291
$ws.Cells.Item(25, 4).Value = 4.907776918959446
$ws.Cells.Item(25, 6).Value = 25.29452162404124
$ws.Cells.Item(25, 7).Value = 3.628407607369999
$ws.Cells.Item(25, 12).Value = 10.82297552751647
$ws.Cells.Item(25, 13).Value = 15.39348601542303
$ws.Cells.Item(25, 14).Value = 18.09747940707215
$ws.Cells.Item(25, 15).Value = 22.4278454343175
